# Actualización desde MV -datos-
# Adds rows 86:106 (dates 06-08-2021 .. 06-09-2021) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Serie (date text), Cupo, Monto demandado, Total adjudicado,
# Adjudicado bancos, Adjudicado AFP, Tasa (or $null when the source cell is blank).
$rows = @(
    @("06-08-2021", 900000,  $null,    0,       $null,  $null,  $null),
    @("09-08-2021", 1000000, 937000,   500000,  318000, 182000, 0.8100000000000001),
    @("10-08-2021", 900000,  906800,   720000,  503200, 216800, 0.86),
    @("11-08-2021", 1000000, 957000,   700000,  573000, 127000, 0.85),
    @("12-08-2021", 900000,  1070000,  450000,  275000, 175000, 0.83),
    @("13-08-2021", 1000000, 944000,   500000,  376000, 124000, 0.92),
    @("16-08-2021", 800000,  920000,   400000,  270000, 130000, 0.87),
    @("17-08-2021", 800000,  1174000,  400000,  286000, 114000, 0.88),
    @("18-08-2021", 800000,  1303000,  800000,  688000, 112000, 0.89),
    @("19-08-2021", 800000,  1069000,  400000,  301000, 99000,  0.9),
    @("20-08-2021", 800000,  1300000,  400000,  350000, 50000,  0.92),
    @("23-08-2021", 800000,  1014000,  800000,  736000, 64000,  1),
    @("24-08-2021", 800000,  1171000,  400000,  200000, 200000, 0.99),
    @("25-08-2021", 800000,  917000,   400000,  333000, 67000,  0.97),
    @("26-08-2021", 900000,  1151000,  900000,  799000, 101000, 0.97),
    @("27-08-2021", 900000,  1025000,  900000,  775000, 125000, 1.02),
    @("30-08-2021", 800000,  $null,    0,       $null,  $null,  $null),
    @("01-09-2021", 900000,  1558000,  1350000, 1157000, 193000, 1.49),
    @("02-09-2021", 900000,  1606000,  1350000, 984000,  366000, 1.5),
    @("03-09-2021", 900000,  1526000,  1350000, 1154000, 196000, 1.49),
    @("06-09-2021", 1000000, 1716000,  1500000, 1174000, 326000, 1.5)
)

$startRow = 86
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A holds the date as plain text (matches the existing "Serie"
    # column), so force text formatting before writing it and then drop
    # back to the default style so no extra style index is introduced.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $data[1]
    if ($null -ne $data[2]) { $ws.Cells.Item($r, 3).Value = $data[2] }
    $ws.Cells.Item($r, 4).Value = $data[3]
    if ($null -ne $data[4]) { $ws.Cells.Item($r, 5).Value = $data[4] }
    if ($null -ne $data[5]) { $ws.Cells.Item($r, 6).Value = $data[5] }
    if ($null -ne $data[6]) { $ws.Cells.Item($r, 7).Value = $data[6] }
}

Write-Output "Added rows 86-106"
